$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 593, shifting existing rows (593:634) down to (594:635)
$ws.Rows.Item(593).Insert()

# Write the new record's values. The date-like text in column A must stay a
# plain text value (matching every other row in the sheet) rather than being
# auto-converted into a date serial number, so it is entered with a leading
# apostrophe (forces text) and the cell style is then reset to Normal so no
# stray "Text" number-format style sticks to the cell.
$ws.Range("A593").Value = "'2026/01/10"
$ws.Range("A593").Style = "Normal"
$ws.Range("B593").Value = "土"
$ws.Range("C593").Value = 12
$ws.Range("D593").Value = 60
